$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.118.79"
$ws.Range("E2").Value = "  -4.78%  "
$ws.Range("D3").Value = "3.402.98"
$ws.Range("E3").Value = "  -3.90%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'560.01"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "'172.19"
$ws.Range("E6").Value = "  -9.41%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("D11").Value = "'56.33"
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("E13").Value = "  -4.02%  "
$ws.Range("D14").Value = "3.947.76"
$ws.Range("E14").Value = "  -3.94%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "'0.119"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.394.01"
$ws.Range("E16").Value = "  -4.41%  "
$ws.Range("D17").Value = "'17.92"
$ws.Range("E17").Value = "  -3.10%  "
$ws.Range("E18").Value = "  -2.50%  "
$ws.Range("D19").Value = "64.108.35"
$ws.Range("E19").Value = "  -4.85%  "
$ws.Range("D20").Value = "'0.984"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "'408.46"
$ws.Range("E21").Value = "  -4.38%  "
$ws.Range("D22").Value = "'4.10"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").Value = "'4.40"
$ws.Range("E23").Value = "  +5.94%  "
$ws.Range("D24").Value = "'13.23"
$ws.Range("E24").Value = "  +7.84%  "
$ws.Range("D25").Value = "'82.65"
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("D26").Value = "'10.74"
$ws.Range("E26").Value = "  -3.28%  "
$ws.Range("D27").Value = "'2.75"
$ws.Range("E27").Value = "  -5.52%  "
$ws.Range("D28").Value = "'8.76"
$ws.Range("E28").Value = "  -2.85%  "
$ws.Range("D29").Value = "'29.52"
$ws.Range("E29").Value = "  -3.31%  "
$ws.Range("D30").Value = "'6.63"
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("D31").Value = "'584.89"
$ws.Range("E31").Value = "  -7.48%  "
$ws.Range("E32").Value = "  -2.60%  "
$ws.Range("E33").Value = "  -4.15%  "
$ws.Range("D34").Value = "'58.78"
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").Value = "'35.94"
$ws.Range("E37").Value = "  -6.56%  "
$ws.Range("E38").Value = "  -4.53%  "
$ws.Range("D39").Value = "3.183.85"
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").Value = "0.0₃0726"
$ws.Range("E41").Value = "  -10.59%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("E44").Value = "  -5.45%  "
$ws.Range("D45").Value = "'3.25"
$ws.Range("E45").Value = "  -3.11%  "
$ws.Range("E46").Value = "  -3.74%  "
$ws.Range("E47").Value = "  -6.17%  "
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'8.23"
$ws.Range("E49").Value = "  -4.61%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'134.03"
$ws.Range("E50").Value = "  -5.13%  "
$ws.Range("E51").Value = "  +2.44%  "
